# Update "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# "zh-cn" and "de-de" worksheets, as part of generating the
# handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 09:51:49"
$wsZhCn.Range("E3").Value = "2016-03-14 09:51:49"
$wsZhCn.Range("H2").Value = "2016-03-14 09:52:40"
$wsZhCn.Range("H3").Value = "2016-03-14 09:52:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 09:51:56"
$wsDeDe.Range("E3").Value = "2016-03-14 09:51:56"
$wsDeDe.Range("H2").Value = "2016-03-14 09:52:54"
$wsDeDe.Range("H3").Value = "2016-03-14 09:52:54"
